$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-06 16:17:46'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '88%'
$ws.Range('K2').Value = '8.6 MJ/m2'
$ws.Range('O2').Value = '-0.1 °C'
$ws.Range('E3').Value = '2026-02-06 16:17:49'
$ws.Range('H3').NumberFormat = '@'
$ws.Range('H3').Value = '70%'
$ws.Range('K3').Value = '12.2 MJ/m2'
$ws.Range('E4').Value = '2026-02-06 16:17:51'
$ws.Range('J4').Value = '996.8 hPa'
$ws.Range('K4').Value = '11.2 MJ/m2'
$ws.Range('O4').Value = '13.6 °C'
$ws.Range('E5').Value = '2026-02-06 16:17:54'
$ws.Range('J5').Value = '997.1 hPa'
$ws.Range('K5').Value = '10.3 MJ/m2'
$ws.Range('O5').Value = '10.9 °C'
$ws.Range('E6').Value = '2026-02-06 16:17:56'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '46%'
$ws.Range('J6').Value = '998.3 hPa'
$ws.Range('K6').Value = '9.7 MJ/m2'
$ws.Range('O6').Value = '15.5 °C'
$ws.Range('E7').Value = '2026-02-06 16:17:59'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '61%'
$ws.Range('J7').Value = '997.9 hPa'
$ws.Range('K7').Value = '11.9 MJ/m2'
$ws.Range('O7').Value = '11.6 °C'
$ws.Range('E8').Value = '2026-02-06 16:18:01'
$ws.Range('K8').Value = '11.7 MJ/m2'
$ws.Range('O8').Value = '10.1 °C'
$ws.Range('E9').Value = '2026-02-06 16:18:04'
$ws.Range('O9').Value = '4.7 °C'
$ws.Range('E10').Value = '2026-02-06 16:18:06'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '87%'
$ws.Range('I10').Value = '0.0 mm'
$ws.Range('M10').Value = '16.1 °C 13:31 TU'
$ws.Range('N10').Value = '2.9 °C 6:09 TU'
$ws.Range('O10').Value = '9.0 °C'
$ws.Range('E11').Value = '2026-02-06 16:18:08'
$ws.Range('K11').Value = '8.8 MJ/m2'
$ws.Range('O11').Value = '5.3 °C'
$ws.Range('E12').Value = '2026-02-06 16:18:11'
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H12').Value = '55%'
$ws.Range('K12').Value = '12.1 MJ/m2'
$ws.Range('E13').Value = '2026-02-06 16:18:13'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value = '75%'
$ws.Range('O13').Value = '10.2 °C'
$ws.Range('E14').Value = '2026-02-06 16:18:16'
$ws.Range('K14').Value = '7.4 MJ/m2'
$ws.Range('E15').Value = '2026-02-06 16:18:18'
$ws.Range('J15').Value = '997.2 hPa'
$ws.Range('K15').Value = '11.5 MJ/m2'
$ws.Range('O15').Value = '10.4 °C'
$ws.Range('E16').Value = '2026-02-06 16:18:21'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '85%'
$ws.Range('K16').Value = '9.5 MJ/m2'
$ws.Range('O16').Value = '5.9 °C'
$ws.Range('E17').Value = '2026-02-06 16:18:24'
$ws.Range('K17').Value = '10.5 MJ/m2'
$ws.Range('O17').Value = '5.9 °C'
$ws.Range('E18').Value = '2026-02-06 16:18:26'
$ws.Range('K18').Value = '5.8 MJ/m2'
$ws.Range('E19').Value = '2026-02-06 16:18:29'
$ws.Range('K19').Value = '11.6 MJ/m2'
$ws.Range('O19').Value = '9.8 °C'
$ws.Range('E20').Value = '2026-02-06 16:18:31'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '77%'
$ws.Range('K20').Value = '12.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-06 16:18:34'
$ws.Range('K21').Value = '10.2 MJ/m2'
$ws.Range('O21').Value = '8.4 °C'
$ws.Range('E22').Value = '2026-02-06 16:18:36'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '75%'
$ws.Range('K22').Value = '11.5 MJ/m2'
$ws.Range('O22').Value = '10.7 °C'
$ws.Range('E23').Value = '2026-02-06 16:18:38'
$ws.Range('J23').Value = '997.2 hPa'
$ws.Range('K23').Value = '8.8 MJ/m2'
$ws.Range('O23').Value = '10.0 °C'
$ws.Range('E24').Value = '2026-02-06 16:18:41'
$ws.Range('J24').Value = '996.6 hPa'
$ws.Range('K24').Value = '11.7 MJ/m2'
$ws.Range('E25').Value = '2026-02-06 16:18:43'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '79%'
$ws.Range('K25').Value = '9.7 MJ/m2'
$ws.Range('L25').Value = '23.4 km/h - 230º 15:35 TU'
$ws.Range('O25').Value = '4.3 °C'
$ws.Range('E26').Value = '2026-02-06 16:18:46'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '79%'
$ws.Range('K26').Value = '7.9 MJ/m2'
$ws.Range('E27').Value = '2026-02-06 16:18:48'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '81%'
$ws.Range('J27').Value = '997.2 hPa'
$ws.Range('K27').Value = '10.4 MJ/m2'
$ws.Range('O27').Value = '10.8 °C'
$ws.Range('E28').Value = '2026-02-06 16:18:51'
$ws.Range('J28').Value = '999.5 hPa'
$ws.Range('O28').Value = '4.7 °C'
$ws.Range('E29').Value = '2026-02-06 16:18:53'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '60%'
$ws.Range('K29').Value = '12.0 MJ/m2'
$ws.Range('O29').Value = '12.7 °C'
$ws.Range('E30').Value = '2026-02-06 16:18:56'
$ws.Range('E31').Value = '2026-02-06 16:18:58'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '84%'
$ws.Range('O31').Value = '7.2 °C'
$ws.Range('E32').Value = '2026-02-06 16:19:01'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '48%'
$ws.Range('J32').Value = '998.6 hPa'
$ws.Range('K32').Value = '12.0 MJ/m2'
$ws.Range('E33').Value = '2026-02-06 16:19:03'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '83%'
$ws.Range('O33').Value = '10.1 °C'
$ws.Range('E34').Value = '2026-02-06 16:19:06'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '75%'
$ws.Range('K34').Value = '11.9 MJ/m2'
$ws.Range('O34').Value = '8.6 °C'
$ws.Range('E35').Value = '2026-02-06 16:19:09'
$ws.Range('K35').Value = '9.3 MJ/m2'
$ws.Range('O35').Value = '-2.0 °C'
$ws.Range('E36').Value = '2026-02-06 16:19:11'
$ws.Range('K36').Value = '11.2 MJ/m2'
$ws.Range('O36').Value = '13.4 °C'
